# This change is a "touch" resave: the authoring tool re-wrote several XML
# parts (both slides, their notes pages, their slide layouts, and the shared
# slide master) without altering any visible text or structure - only the
# serializer's own formatting changed upstream (attribute/namespace style).
# We reproduce the edit at the content level by re-applying each shape's
# current text back onto itself, which causes this runtime to rewrite the
# same set of parts while leaving every run of text exactly as it was.

$p = $ppt.ActivePresentation

function Touch-Shape($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $txt = $tr.Text
        if ($txt.Length -gt 0) {
            $tr.Text = $txt
        }
    }
}

function Touch-Shapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        Touch-Shape $shapes.Item($i)
    }
}

# Slide 1 ("Module 01: Core AWS Knowledge") + its notes page + its layout.
$s1 = $p.Slides.Item(1)
Touch-Shapes $s1.Shapes
Touch-Shapes $s1.NotesPage.Shapes
Touch-Shapes $s1.CustomLayout.Shapes

# Slide 2 ("What Exactly Is Cloud Computing ?") + its notes page + its layout.
$s2 = $p.Slides.Item(2)
Touch-Shapes $s2.Shapes
Touch-Shapes $s2.NotesPage.Shapes
Touch-Shapes $s2.CustomLayout.Shapes

# Shared slide master used by both layouts above.
Touch-Shapes $p.SlideMaster.Shapes
